# Remove every <w:contextualSpacing w:val="0"/> element from the
# paragraph properties (w:pPr) throughout the document - this affects
# both the main document story (word/document.xml) and the comments
# part (word/comments.xml).
#
# The Word object model exposed by this host does not surface a
# ParagraphFormat.ContextualSpacing property, so we go through the OOXML
# round-trip that *is* exposed: Range.WordOpenXML (read) / Range.InsertXML
# (write). Reading/writing through Document.Content is the only Range
# whose InsertXML call is actually wired to every package part (including
# word/comments.xml); scoping the same call to a single paragraph or a
# single Comment.Range only ever touches word/document.xml (or nothing at
# all), so we rebuild a minimal OOXML package containing just the parts
# that need the edit and push that back through Document.Content.

$d = $word.ActiveDocument
$full = $d.Content.WordOpenXML

function Get-OpcPart($packageXml, $partName) {
    $marker = '<pkg:part pkg:name="' + $partName + '"'
    $start = $packageXml.IndexOf($marker)
    if ($start -lt 0) { return $null }
    $closeTag = "</pkg:part>"
    $end = $packageXml.IndexOf($closeTag, $start)
    return $packageXml.Substring($start, $end - $start + $closeTag.Length)
}

function Remove-ContextualSpacing($partXml) {
    if ($partXml -eq $null) { return $null }
    return $partXml.Replace('<w:contextualSpacing w:val="0"/>', "")
}

$docPart = Get-OpcPart $full "/word/document.xml"
$commentsPart = Get-OpcPart $full "/word/comments.xml"

$docPart = Remove-ContextualSpacing $docPart
$commentsPart = Remove-ContextualSpacing $commentsPart

$rebuilt = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + $docPart

if ($commentsPart -ne $null) {
    $rebuilt = $rebuilt + $commentsPart
}

$rebuilt = $rebuilt + '</pkg:package>'

$d.Content.InsertXML($rebuilt)
